# Case_1_153 / res_bus / vm_pu: the 380 kV slack-bus setpoint was lowered from
# 1.05 p.u. to 1.02 p.u. (column B), which changes the converged load-flow solution
# for every other bus voltage magnitude (columns C-F and I-N) across all time steps
# (rows 2-25). Column G (=1, reference) and the header/time-step columns are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (time step 0)
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.053854969768837
$ws.Cells.Item(2, 4).Value = 1.065385528528399
$ws.Cells.Item(2, 5).Value = 1.050444478685133
$ws.Cells.Item(2, 6).Value = 1.072625718102157
$ws.Cells.Item(2, 9).Value = 1.05070393706763
$ws.Cells.Item(2, 10).Value = 1.058870125315519
$ws.Cells.Item(2, 11).Value = 1.068098910445472
$ws.Cells.Item(2, 12).Value = 1.053198738339188
$ws.Cells.Item(2, 13).Value = 1.075319734315271
$ws.Cells.Item(2, 14).Value = 1.060373842347568

# Row 3 (time step 1)
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.055213145730025
$ws.Cells.Item(3, 4).Value = 1.066268974430999
$ws.Cells.Item(3, 5).Value = 1.051615880079835
$ws.Cells.Item(3, 6).Value = 1.073703333775923
$ws.Cells.Item(3, 9).Value = 1.051090693577514
$ws.Cells.Item(3, 10).Value = 1.05987736903089
$ws.Cells.Item(3, 11).Value = 1.068797574889382
$ws.Cells.Item(3, 12).Value = 1.054181662527473
$ws.Cells.Item(3, 13).Value = 1.076213487385941
$ws.Cells.Item(3, 14).Value = 1.061382516464547

# Row 4 (time step 2)
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.056091382187822
$ws.Cells.Item(4, 4).Value = 1.066840266128102
$ws.Cells.Item(4, 5).Value = 1.052373531253224
$ws.Cells.Item(4, 6).Value = 1.074400505318449
$ws.Cells.Item(4, 9).Value = 1.051339513624125
$ws.Cells.Item(4, 10).Value = 1.060528057922242
$ws.Cells.Item(4, 11).Value = 1.069248661828744
$ws.Cells.Item(4, 12).Value = 1.054816782168625
$ws.Cells.Item(4, 13).Value = 1.076791086608687
$ws.Cells.Item(4, 14).Value = 1.06203412940876

# Row 5 (time step 3)
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.056460454426258
$ws.Cells.Item(5, 4).Value = 1.067080352540554
$ws.Cells.Item(5, 5).Value = 1.052691972817789
$ws.Cells.Item(5, 6).Value = 1.074693569394882
$ws.Cells.Item(5, 9).Value = 1.05144377432063
$ws.Cells.Item(5, 10).Value = 1.060801355072274
$ws.Cells.Item(5, 11).Value = 1.06943806137775
$ws.Cells.Item(5, 12).Value = 1.055083573813426
$ws.Cells.Item(5, 13).Value = 1.077033738584711
$ws.Cells.Item(5, 14).Value = 1.062307814672093

# Row 6 (time step 4)
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.056522415277903
$ws.Cells.Item(6, 4).Value = 1.067120659167354
$ws.Cells.Item(6, 5).Value = 1.05274543623411
$ws.Cells.Item(6, 6).Value = 1.074742774608786
$ws.Cells.Item(6, 9).Value = 1.051461260046177
$ws.Cells.Item(6, 10).Value = 1.060847228135595
$ws.Cells.Item(6, 11).Value = 1.069469848502136
$ws.Cells.Item(6, 12).Value = 1.05512835690235
$ws.Cells.Item(6, 13).Value = 1.077074470945911
$ws.Cells.Item(6, 14).Value = 1.062353752880425

# Row 7 (time step 5)
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.056096314287989
$ws.Cells.Item(7, 4).Value = 1.066843474505972
$ws.Cells.Item(7, 5).Value = 1.052377786576487
$ws.Cells.Item(7, 6).Value = 1.074404421358739
$ws.Cells.Item(7, 9).Value = 1.05134090810836
$ws.Cells.Item(7, 10).Value = 1.060531710720979
$ws.Cells.Item(7, 11).Value = 1.069251193526413
$ws.Cells.Item(7, 12).Value = 1.054820347884314
$ws.Cells.Item(7, 13).Value = 1.076794329605487
$ws.Cells.Item(7, 14).Value = 1.06203778739489

# Row 8 (time step 6)
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.0543140970699
$ws.Cells.Item(8, 4).Value = 1.065684167597087
$ws.Cells.Item(8, 5).Value = 1.050840427861336
$ws.Cells.Item(8, 6).Value = 1.072989928597563
$ws.Cells.Item(8, 9).Value = 1.050834941342777
$ws.Cells.Item(8, 10).Value = 1.059210750120145
$ws.Cells.Item(8, 11).Value = 1.068335234109851
$ws.Cells.Item(8, 12).Value = 1.053531109524756
$ws.Cells.Item(8, 13).Value = 1.07562193198622
$ws.Cells.Item(8, 14).Value = 1.060714950878487

# Row 9 (time step 7)
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.051168877443448
$ws.Cells.Item(9, 4).Value = 1.063638549993999
$ws.Cells.Item(9, 5).Value = 1.048128810662596
$ws.Cells.Item(9, 6).Value = 1.070496446728738
$ws.Cells.Item(9, 9).Value = 1.049932317708758
$ws.Cells.Item(9, 10).Value = 1.056874776987457
$ws.Cells.Item(9, 11).Value = 1.066713527087278
$ws.Cells.Item(9, 12).Value = 1.051252323791127
$ws.Cells.Item(9, 13).Value = 1.073550452023513
$ws.Cells.Item(9, 14).Value = 1.058375660396011

# Row 10 (time step 8)
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.049068641987648
$ws.Cells.Item(10, 4).Value = 1.062272886808031
$ws.Cells.Item(10, 5).Value = 1.046319159200043
$ws.Cells.Item(10, 6).Value = 1.068833386807076
$ws.Cells.Item(10, 9).Value = 1.049323088141541
$ws.Cells.Item(10, 10).Value = 1.055311738655989
$ws.Cells.Item(10, 11).Value = 1.065627168982403
$ws.Cells.Item(10, 12).Value = 1.049728289300288
$ws.Cells.Item(10, 13).Value = 1.072165636180277
$ws.Cells.Item(10, 14).Value = 1.056810402370829

# Row 11 (time step 9)
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.048158349939417
$ws.Cells.Item(11, 4).Value = 1.061681071434361
$ws.Cells.Item(11, 5).Value = 1.04553507113331
$ws.Cells.Item(11, 6).Value = 1.068113067390117
$ws.Cells.Item(11, 9).Value = 1.049057498280904
$ws.Cells.Item(11, 10).Value = 1.054633531860031
$ws.Cells.Item(11, 11).Value = 1.065155510938441
$ws.Cells.Item(11, 12).Value = 1.049067182823979
$ws.Cells.Item(11, 13).Value = 1.07156506719833
$ws.Cells.Item(11, 14).Value = 1.056131232443429

# Row 12 (time step 10)
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.047820090460495
$ws.Cells.Item(12, 4).Value = 1.061461172372349
$ws.Cells.Item(12, 5).Value = 1.045243747981614
$ws.Cells.Item(12, 6).Value = 1.067845476208698
$ws.Cells.Item(12, 9).Value = 1.048958576398692
$ws.Cells.Item(12, 10).Value = 1.054381401953852
$ws.Cells.Item(12, 11).Value = 1.064980125796466
$ws.Cells.Item(12, 12).Value = 1.048821436738487
$ws.Cells.Item(12, 13).Value = 1.071341846889252
$ws.Cells.Item(12, 14).Value = 1.055878744483864

# Row 13 (time step 11)
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.047892654520919
$ws.Cells.Item(13, 4).Value = 1.06150834476757
$ws.Cells.Item(13, 5).Value = 1.04530624136564
$ws.Cells.Item(13, 6).Value = 1.067902876925884
$ws.Cells.Item(13, 9).Value = 1.048979807713867
$ws.Cells.Item(13, 10).Value = 1.054435494397354
$ws.Cells.Item(13, 11).Value = 1.065017755140363
$ws.Cells.Item(13, 12).Value = 1.048874158378823
$ws.Cells.Item(13, 13).Value = 1.071389734869992
$ws.Cells.Item(13, 14).Value = 1.05593291374484

# Row 14 (time step 12)
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.048130392117595
$ws.Cells.Item(14, 4).Value = 1.061662895995105
$ws.Cells.Item(14, 5).Value = 1.045510991879402
$ws.Cells.Item(14, 6).Value = 1.068090948893185
$ws.Cells.Item(14, 9).Value = 1.049049326881437
$ws.Cells.Item(14, 10).Value = 1.054612695101634
$ws.Cells.Item(14, 11).Value = 1.065141017438817
$ws.Cells.Item(14, 12).Value = 1.049046873099997
$ws.Cells.Item(14, 13).Value = 1.071546618650456
$ws.Cells.Item(14, 14).Value = 1.056110366094445

# Row 15 (time step 13)
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.048276851867102
$ws.Cells.Item(15, 4).Value = 1.061758110465463
$ws.Cells.Item(15, 5).Value = 1.045637135004072
$ws.Cells.Item(15, 6).Value = 1.068206821878584
$ws.Cells.Item(15, 9).Value = 1.049092124125498
$ws.Cells.Item(15, 10).Value = 1.054721845888801
$ws.Cells.Item(15, 11).Value = 1.065216938145469
$ws.Cells.Item(15, 12).Value = 1.049153264172502
$ws.Cells.Item(15, 13).Value = 1.071643261020876
$ws.Cells.Item(15, 14).Value = 1.05621967188825

# Row 16 (time step 14)
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.049129036199212
$ws.Cells.Item(16, 4).Value = 1.062312153548978
$ws.Cells.Item(16, 5).Value = 1.046371185792937
$ws.Cells.Item(16, 6).Value = 1.068881187584751
$ws.Cells.Item(16, 9).Value = 1.049340676667334
$ws.Cells.Item(16, 10).Value = 1.055356719216354
$ws.Cells.Item(16, 11).Value = 1.065658444755857
$ws.Cells.Item(16, 12).Value = 1.049772139443187
$ws.Cells.Item(16, 13).Value = 1.072205474118447
$ws.Cells.Item(16, 14).Value = 1.056855446808749

# Row 17 (time step 15)
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.049663350913751
$ws.Cells.Item(17, 4).Value = 1.062659562309366
$ws.Cells.Item(17, 5).Value = 1.046831501092093
$ws.Cells.Item(17, 6).Value = 1.06930414355282
$ws.Cells.Item(17, 9).Value = 1.049496107223655
$ws.Cells.Item(17, 10).Value = 1.055754581143331
$ws.Cells.Item(17, 11).Value = 1.065935052469657
$ws.Cells.Item(17, 12).Value = 1.050160022978745
$ws.Cells.Item(17, 13).Value = 1.072557884031992
$ws.Cells.Item(17, 14).Value = 1.057253873745298

# Row 18 (time step 16)
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.049974923431394
$ws.Cells.Item(18, 4).Value = 1.062862154293323
$ws.Cells.Item(18, 5).Value = 1.047099947482718
$ws.Cells.Item(18, 6).Value = 1.069550827109209
$ws.Cells.Item(18, 9).Value = 1.049586594660489
$ws.Cells.Item(18, 10).Value = 1.055986512329428
$ws.Cells.Item(18, 11).Value = 1.066096271802321
$ws.Cells.Item(18, 12).Value = 1.05038615431865
$ws.Cells.Item(18, 13).Value = 1.072763348724302
$ws.Cells.Item(18, 14).Value = 1.057486134300282

# Row 19 (time step 17)
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.05008114746938
$ws.Cells.Item(19, 4).Value = 1.062931225215128
$ws.Cells.Item(19, 5).Value = 1.047191472765771
$ws.Cells.Item(19, 6).Value = 1.069634936575821
$ws.Cells.Item(19, 9).Value = 1.049617419300179
$ws.Cells.Item(19, 10).Value = 1.056065572078664
$ws.Cells.Item(19, 11).Value = 1.066151222900491
$ws.Cells.Item(19, 12).Value = 1.050463239919883
$ws.Cells.Item(19, 13).Value = 1.072833391656034
$ws.Cells.Item(19, 14).Value = 1.057565306323431

# Row 20 (time step 18)
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.04960603274888
$ws.Cells.Item(20, 4).Value = 1.062622293370693
$ws.Cells.Item(20, 5).Value = 1.046782118579153
$ws.Cells.Item(20, 6).Value = 1.069258766412439
$ws.Cells.Item(20, 9).Value = 1.049479448855158
$ws.Cells.Item(20, 10).Value = 1.055711908306424
$ws.Cells.Item(20, 11).Value = 1.065905387633637
$ws.Cells.Item(20, 12).Value = 1.050118418618941
$ws.Cells.Item(20, 13).Value = 1.072520083111741
$ws.Cells.Item(20, 14).Value = 1.057211140308068

# Row 21 (time step 19)
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.048060388127325
$ws.Cells.Item(21, 4).Value = 1.061617386521534
$ws.Cells.Item(21, 5).Value = 1.045450700133919
$ws.Cells.Item(21, 6).Value = 1.068035567298291
$ws.Cells.Item(21, 9).Value = 1.049028862681328
$ws.Cells.Item(21, 10).Value = 1.054560519836216
$ws.Cells.Item(21, 11).Value = 1.065104725029638
$ws.Cells.Item(21, 12).Value = 1.048996017959489
$ws.Cells.Item(21, 13).Value = 1.07150042422892
$ws.Cells.Item(21, 14).Value = 1.056058116734166

# Row 22 (time step 20)
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.04708778795864
$ws.Cells.Item(22, 4).Value = 1.060985141722746
$ws.Cells.Item(22, 5).Value = 1.044613131997628
$ws.Cells.Item(22, 6).Value = 1.06726630460784
$ws.Cells.Item(22, 9).Value = 1.048743998665833
$ws.Cells.Item(22, 10).Value = 1.053835358161305
$ws.Cells.Item(22, 11).Value = 1.064600214715352
$ws.Cells.Item(22, 12).Value = 1.048289267336602
$ws.Cells.Item(22, 13).Value = 1.070858500292446
$ws.Cells.Item(22, 14).Value = 1.055331925246499

# Row 23 (time step 21)
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.047603458406402
$ws.Cells.Item(23, 4).Value = 1.061320346877981
$ws.Cells.Item(23, 5).Value = 1.045057186698163
$ws.Cells.Item(23, 6).Value = 1.067674123855245
$ws.Cells.Item(23, 9).Value = 1.048895158968718
$ws.Cells.Item(23, 10).Value = 1.054219898522299
$ws.Cells.Item(23, 11).Value = 1.064867770091235
$ws.Cells.Item(23, 12).Value = 1.048664030009252
$ws.Cells.Item(23, 13).Value = 1.071198875031695
$ws.Cells.Item(23, 14).Value = 1.055717011698913

# Row 24 (time step 22)
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.049631932620328
$ws.Cells.Item(24, 4).Value = 1.062639133741045
$ws.Cells.Item(24, 5).Value = 1.046804432558504
$ws.Cells.Item(24, 6).Value = 1.06927927044925
$ws.Cells.Item(24, 9).Value = 1.049486976588053
$ws.Cells.Item(24, 10).Value = 1.055731190741424
$ws.Cells.Item(24, 11).Value = 1.065918792270851
$ws.Cells.Item(24, 12).Value = 1.050137218192178
$ws.Cells.Item(24, 13).Value = 1.072537163998731
$ws.Cells.Item(24, 14).Value = 1.057230450126337

# Row 25 (time step 23)
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.051982577281191
$ws.Cells.Item(25, 4).Value = 1.064167725169602
$ws.Cells.Item(25, 5).Value = 1.048830152838565
$ws.Cells.Item(25, 6).Value = 1.071141195123336
$ws.Cells.Item(25, 9).Value = 1.050166982058338
$ws.Cells.Item(25, 10).Value = 1.05747967886124
$ws.Cells.Item(25, 11).Value = 1.06713369261127
$ws.Cells.Item(25, 12).Value = 1.051842287349917
$ws.Cells.Item(25, 13).Value = 1.074086647738439
$ws.Cells.Item(25, 14).Value = 1.058981421299837
